$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 22.666666
$ws.Range("I11").Value = 22.666666
$ws.Range("K11").Value = 22.666666
$ws.Range("M11").Value = 117.333334

$ws.Range("H34").Value = 2797.8
$ws.Range("I34").Value = 2797.8
$ws.Range("K34").Value = 2797.8
$ws.Range("M34").Value = -2594.8

$ws.Range("H36").Value = 2797.8
$ws.Range("I36").Value = 2797.8
$ws.Range("K36").Value = 2797.8
$ws.Range("M36").Value = -2082.8

$ws.Range("H47").Value = 2067
$ws.Range("I47").Value = 2067
$ws.Range("K47").Value = 2067
$ws.Range("M47").Value = -1095

$ws.Range("H64").Value = 6874.375
$ws.Range("I64").Value = 4999
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 4999
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -4751
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 6874.375
$ws.Range("I67").Value = 4999
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 4999
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -4141
$ws.Range("N67").Value = -11716

$ws.Range("H88").Value = 1438.6666
$ws.Range("J88").Value = 1513.375
$ws.Range("L88").Value = 1513.375
$ws.Range("N88").Value = -2325.375

$ws.Range("H91").Value = 1438.6666
$ws.Range("J91").Value = 1513.375
$ws.Range("L91").Value = 1513.375
$ws.Range("N91").Value = -4321.375

$ws.Range("H94").Value = 3124
$ws.Range("I94").Value = 3124
$ws.Range("K94").Value = 3124
$ws.Range("M94").Value = -2673

$ws.Range("H100").Value = 6000.25
$ws.Range("I100").Value = 4000.5
$ws.Range("K100").Value = 4000.5
$ws.Range("M100").Value = -3459.5

$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254

$ws.Range("H125").Value = 1676.75
$ws.Range("J125").Value = 2537.5
$ws.Range("L125").Value = 22837.5
$ws.Range("N125").Value = -27757.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 999
$ws.Range("I4").Value = 998
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 998
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -882
$ws.Range("N4").Value = -1232

$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88

$ws.Range("H61").Value = 3954.4546
$ws.Range("I61").Value = 3954.4546
$ws.Range("K61").Value = 3954.4546
$ws.Range("M61").Value = -3742.4546

$ws.Range("H98").Value = 8737.5
$ws.Range("J98").Value = 8737.5
$ws.Range("L98").Value = 8737.5
$ws.Range("N98").Value = -14727.5

$ws.Range("H122").Value = 2065.25
$ws.Range("I122").Value = 1001.4
$ws.Range("K122").Value = 3004.2
$ws.Range("M122").Value = -554.1999999999998

$ws.Range("H136").Value = 3954.4546
$ws.Range("I136").Value = 3954.4546
$ws.Range("K136").Value = 11863.3638
$ws.Range("M136").Value = -9313.363799999999

$ws.Range("H139").Value = 45000
$ws.Range("J139").Value = 45000
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0

$ws.Range("H74").Value = 145000
$ws.Range("J74").Value = 145000
$ws.Range("L74").Value = 145000
$ws.Range("N74").Value = -146872

$ws.Range("H77").Value = 145000
$ws.Range("J77").Value = 145000
$ws.Range("L77").Value = 435000
$ws.Range("N77").Value = -444360

$ws.Range("H81").Value = 93600
$ws.Range("J81").Value = 93600
$ws.Range("L81").Value = 93600
$ws.Range("N81").Value = -95722

$ws.Range("H84").Value = 93600
$ws.Range("J84").Value = 93600
$ws.Range("L84").Value = 280800
$ws.Range("N84").Value = -291408

$ws.Range("H107").Value = 5748.357
$ws.Range("I107").Value = 1622.25
$ws.Range("J107").Value = 7398.8
$ws.Range("K107").Value = 1622.25
$ws.Range("L107").Value = 7398.8
$ws.Range("M107").Value = 297.75
$ws.Range("N107").Value = -11238.8

$ws.Range("H134").Value = 1463.3334
$ws.Range("I134").Value = 1200
$ws.Range("J134").Value = 1990
$ws.Range("K134").Value = 3600
$ws.Range("L134").Value = 5970
$ws.Range("M134").Value = -1065
$ws.Range("N134").Value = -11040

$ws.Range("H137").Value = 135520
$ws.Range("J137").Value = 135520
$ws.Range("L137").Value = 135520
$ws.Range("N137").Value = -145720

$ws.Range("H138").Value = 200000
$ws.Range("J138").Value = 200000
$ws.Range("L138").Value = 200000
$ws.Range("N138").Value = -210280

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws.Range("H140").Value = 145390
$ws.Range("J140").Value = 145390
$ws.Range("L140").Value = 145390
$ws.Range("N140").Value = -155750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4453.3335
$ws.Range("J58").Value = 7520
$ws.Range("L58").Value = 7520
$ws.Range("N58").Value = -7926

$ws.Range("H132").Value = 4300
$ws.Range("I132").Value = 4300
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12900
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -10370

$ws.Range("H136").Value = 4453.3335
$ws.Range("J136").Value = 7520
$ws.Range("L136").Value = 22560
$ws.Range("N136").Value = -27660

$ws.Range("H141").Value = 80984.7
$ws.Range("J141").Value = 80984.7
$ws.Range("L141").Value = 80984.7
$ws.Range("N141").Value = -91344.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 7500
$ws.Range("I12").Value = 7000
$ws.Range("J12").Value = 8750
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 8750
$ws.Range("M12").Value = -6860
$ws.Range("N12").Value = -9030

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("N18").Value = 0

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H97").Value = 1209.4
$ws.Range("I97").Value = 950
$ws.Range("J97").Value = 1274.25
$ws.Range("K97").Value = 950
$ws.Range("L97").Value = 1274.25
$ws.Range("M97").Value = -454
$ws.Range("N97").Value = -2266.25

$ws.Range("H126").Value = 3561.3125
$ws.Range("I126").Value = 2690.3333
$ws.Range("K126").Value = 8070.999899999999
$ws.Range("M126").Value = -5600.999899999999

$ws.Range("H132").Value = 4725.2
$ws.Range("I132").Value = 4855.125
$ws.Range("J132").Value = 4205.5
$ws.Range("K132").Value = 14565.375
$ws.Range("L132").Value = 12616.5
$ws.Range("M132").Value = -12035.375
$ws.Range("N132").Value = -17676.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("N2").Value = 0

$ws.Range("H7").Value = 8997.5
$ws.Range("I7").Value = 7999.5
$ws.Range("J7").Value = 9995.5
$ws.Range("K7").Value = 7999.5
$ws.Range("L7").Value = 9995.5
$ws.Range("M7").Value = -7887.5
$ws.Range("N7").Value = -10219.5

$ws.Range("H126").Value = 8997.5
$ws.Range("I126").Value = 7999.5
$ws.Range("J126").Value = 9995.5
$ws.Range("K126").Value = 23998.5
$ws.Range("L126").Value = 29986.5
$ws.Range("M126").Value = -21528.5
$ws.Range("N126").Value = -34926.5

$ws.Range("H132").Value = 5198
$ws.Range("I132").Value = 5497.75
$ws.Range("K132").Value = 16493.25
$ws.Range("M132").Value = -13963.25

$ws.Range("H136").Value = 4
$ws.Range("I136").Value = 4
$ws.Range("K136").Value = 12
$ws.Range("M136").Value = 2538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 946.25
$ws.Range("I2").Value = 946.25
$ws.Range("K2").Value = 946.25
$ws.Range("M2").Value = -834.25
